# Edit script for timings_inference.xlsx
# Adds a new timing value (H56) with an empty styled time cell block (H57:H58),
# and appends new benchmark rows for the "pannuke"/"hovernet" dataset/model
# combination (including some OOM entries), per "fixed hp search for other
# metrics".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New time-formatted cell block (H56:H58) ---
$ws.Range("H56").Value2 = 0.56680555555555556
$ws.Range("H56").NumberFormat = "h:mm:ss"
$ws.Range("H57").NumberFormat = "h:mm:ss"
$ws.Range("H58").NumberFormat = "h:mm:ss"

# --- New benchmark data rows (91-100, 103-112) ---
# Row 91
$ws.Cells.Item(91,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(91,2).Value2 = "RTX3090"
$ws.Cells.Item(91,3).Value2 = 4
$ws.Cells.Item(91,4).Value2 = 32
$ws.Cells.Item(91,5).Value2 = "pannuke"
$ws.Cells.Item(91,6).Value2 = "tiny"
$ws.Cells.Item(91,7).Value2 = 4
$ws.Cells.Item(91,8).Value2 = "inference"
$ws.Cells.Item(91,9).Value2 = "TCGA-AA-3977-01Z-00-DX1"
$ws.Cells.Item(91,10).Value2 = 58
$ws.Cells.Item(91,11).Value2 = 4
$ws.Cells.Item(91,12).Value2 = "N/A"

# Row 92
$ws.Cells.Item(92,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(92,2).Value2 = "RTX3090"
$ws.Cells.Item(92,3).Value2 = 4
$ws.Cells.Item(92,4).Value2 = 32
$ws.Cells.Item(92,5).Value2 = "pannuke"
$ws.Cells.Item(92,6).Value2 = "tiny"
$ws.Cells.Item(92,7).Value2 = 4
$ws.Cells.Item(92,8).Value2 = "inference"
$ws.Cells.Item(92,9).Value2 = "TCGA-AA-3688-01Z-00-DX1"
$ws.Cells.Item(92,10).Value2 = 127
$ws.Cells.Item(92,11).Value2 = 4
$ws.Cells.Item(92,12).Value2 = "N/A"

# Row 93
$ws.Cells.Item(93,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(93,2).Value2 = "RTX3090"
$ws.Cells.Item(93,3).Value2 = 4
$ws.Cells.Item(93,4).Value2 = 32
$ws.Cells.Item(93,5).Value2 = "pannuke"
$ws.Cells.Item(93,6).Value2 = "tiny"
$ws.Cells.Item(93,7).Value2 = 4
$ws.Cells.Item(93,8).Value2 = "inference"
$ws.Cells.Item(93,9).Value2 = "TCGA-AA-A010-01Z-00-DX1"
$ws.Cells.Item(93,10).Value2 = 235
$ws.Cells.Item(93,11).Value2 = 4
$ws.Cells.Item(93,12).Value2 = "N/A"

# Row 94
$ws.Cells.Item(94,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(94,2).Value2 = "RTX3090"
$ws.Cells.Item(94,3).Value2 = 4
$ws.Cells.Item(94,4).Value2 = 32
$ws.Cells.Item(94,5).Value2 = "pannuke"
$ws.Cells.Item(94,6).Value2 = "tiny"
$ws.Cells.Item(94,7).Value2 = 4
$ws.Cells.Item(94,8).Value2 = "inference"
$ws.Cells.Item(94,9).Value2 = "TCGA-CK-4951-01Z-00-DX1"
$ws.Cells.Item(94,10).Value2 = 391
$ws.Cells.Item(94,11).Value2 = 4
$ws.Cells.Item(94,12).Value2 = "N/A"

# Row 95
$ws.Cells.Item(95,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(95,2).Value2 = "RTX3090"
$ws.Cells.Item(95,3).Value2 = 4
$ws.Cells.Item(95,4).Value2 = 32
$ws.Cells.Item(95,5).Value2 = "pannuke"
$ws.Cells.Item(95,6).Value2 = "tiny"
$ws.Cells.Item(95,7).Value2 = 4
$ws.Cells.Item(95,8).Value2 = "inference"
$ws.Cells.Item(95,9).Value2 = "TCGA-5M-AAT5-01Z-00-DX1"
$ws.Cells.Item(95,10).Value2 = 1017
$ws.Cells.Item(95,11).Value2 = 4
$ws.Cells.Item(95,12).Value2 = "N/A"

# Row 96
$ws.Cells.Item(96,1).Value2 = "XEON E5-2630 v4"
$ws.Cells.Item(96,2).Value2 = "N/A"
$ws.Cells.Item(96,3).Value2 = 20
$ws.Cells.Item(96,4).Value2 = 128
$ws.Cells.Item(96,5).Value2 = "pannuke"
$ws.Cells.Item(96,6).Value2 = "N/A"
$ws.Cells.Item(96,7).Value2 = "N/A"
$ws.Cells.Item(96,8).Value2 = "postproc"
$ws.Cells.Item(96,9).Value2 = "TCGA-AA-3977-01Z-00-DX1"
$ws.Cells.Item(96,10).Value2 = 56
$ws.Cells.Item(96,11).Value2 = 20
$ws.Cells.Item(96,12).Value2 = 10

# Row 97
$ws.Cells.Item(97,1).Value2 = "XEON E5-2630 v5"
$ws.Cells.Item(97,2).Value2 = "N/A"
$ws.Cells.Item(97,3).Value2 = 20
$ws.Cells.Item(97,4).Value2 = 128
$ws.Cells.Item(97,5).Value2 = "pannuke"
$ws.Cells.Item(97,6).Value2 = "N/A"
$ws.Cells.Item(97,7).Value2 = "N/A"
$ws.Cells.Item(97,8).Value2 = "postproc"
$ws.Cells.Item(97,9).Value2 = "TCGA-AA-3688-01Z-00-DX1"
$ws.Cells.Item(97,10).Value2 = 108
$ws.Cells.Item(97,11).Value2 = 20
$ws.Cells.Item(97,12).Value2 = 10

# Row 98
$ws.Cells.Item(98,1).Value2 = "XEON E5-2630 v6"
$ws.Cells.Item(98,2).Value2 = "N/A"
$ws.Cells.Item(98,3).Value2 = 20
$ws.Cells.Item(98,4).Value2 = 128
$ws.Cells.Item(98,5).Value2 = "pannuke"
$ws.Cells.Item(98,6).Value2 = "N/A"
$ws.Cells.Item(98,7).Value2 = "N/A"
$ws.Cells.Item(98,8).Value2 = "postproc"
$ws.Cells.Item(98,9).Value2 = "TCGA-AA-A010-01Z-00-DX1"
$ws.Cells.Item(98,10).Value2 = 200
$ws.Cells.Item(98,11).Value2 = 20
$ws.Cells.Item(98,12).Value2 = 10

# Row 99
$ws.Cells.Item(99,1).Value2 = "XEON E5-2630 v7"
$ws.Cells.Item(99,2).Value2 = "N/A"
$ws.Cells.Item(99,3).Value2 = 20
$ws.Cells.Item(99,4).Value2 = 128
$ws.Cells.Item(99,5).Value2 = "pannuke"
$ws.Cells.Item(99,6).Value2 = "N/A"
$ws.Cells.Item(99,7).Value2 = "N/A"
$ws.Cells.Item(99,8).Value2 = "postproc"
$ws.Cells.Item(99,9).Value2 = "TCGA-CK-4951-01Z-00-DX1"
$ws.Cells.Item(99,10).Value2 = 263
$ws.Cells.Item(99,11).Value2 = 20
$ws.Cells.Item(99,12).Value2 = 10

# Row 100
$ws.Cells.Item(100,1).Value2 = "XEON E5-2630 v8"
$ws.Cells.Item(100,2).Value2 = "N/A"
$ws.Cells.Item(100,3).Value2 = 20
$ws.Cells.Item(100,4).Value2 = 128
$ws.Cells.Item(100,5).Value2 = "pannuke"
$ws.Cells.Item(100,6).Value2 = "N/A"
$ws.Cells.Item(100,7).Value2 = "N/A"
$ws.Cells.Item(100,8).Value2 = "postproc"
$ws.Cells.Item(100,9).Value2 = "TCGA-5M-AAT5-01Z-00-DX1"
$ws.Cells.Item(100,10).Value2 = 542
$ws.Cells.Item(100,11).Value2 = 20
$ws.Cells.Item(100,12).Value2 = 10

# Row 103
$ws.Cells.Item(103,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(103,2).Value2 = "RTX3090"
$ws.Cells.Item(103,3).Value2 = 4
$ws.Cells.Item(103,4).Value2 = 128
$ws.Cells.Item(103,5).Value2 = "pannuke"
$ws.Cells.Item(103,6).Value2 = "hovernet"
$ws.Cells.Item(103,7).Value2 = "N/A"
$ws.Cells.Item(103,8).Value2 = "inference"
$ws.Cells.Item(103,9).Value2 = "TCGA-AA-3977-01Z-00-DX1"
$ws.Cells.Item(103,10).Value2 = 6643

# Row 104
$ws.Cells.Item(104,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(104,2).Value2 = "RTX3090"
$ws.Cells.Item(104,3).Value2 = 4
$ws.Cells.Item(104,4).Value2 = 128
$ws.Cells.Item(104,5).Value2 = "pannuke"
$ws.Cells.Item(104,6).Value2 = "hovernet"
$ws.Cells.Item(104,7).Value2 = "N/A"
$ws.Cells.Item(104,8).Value2 = "inference"
$ws.Cells.Item(104,9).Value2 = "TCGA-AA-3688-01Z-00-DX1"
$ws.Cells.Item(104,10).Value2 = "OOM"

# Row 105
$ws.Cells.Item(105,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(105,2).Value2 = "RTX3090"
$ws.Cells.Item(105,3).Value2 = 4
$ws.Cells.Item(105,4).Value2 = 128
$ws.Cells.Item(105,5).Value2 = "pannuke"
$ws.Cells.Item(105,6).Value2 = "hovernet"
$ws.Cells.Item(105,7).Value2 = "N/A"
$ws.Cells.Item(105,8).Value2 = "inference"
$ws.Cells.Item(105,9).Value2 = "TCGA-AA-A010-01Z-00-DX1"
$ws.Cells.Item(105,10).Value2 = "OOM"

# Row 106
$ws.Cells.Item(106,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(106,2).Value2 = "RTX3090"
$ws.Cells.Item(106,3).Value2 = 4
$ws.Cells.Item(106,4).Value2 = 128
$ws.Cells.Item(106,5).Value2 = "pannuke"
$ws.Cells.Item(106,6).Value2 = "hovernet"
$ws.Cells.Item(106,7).Value2 = "N/A"
$ws.Cells.Item(106,8).Value2 = "inference"
$ws.Cells.Item(106,9).Value2 = "TCGA-CK-4951-01Z-00-DX1"
$ws.Cells.Item(106,10).Value2 = "OOM"

# Row 107
$ws.Cells.Item(107,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(107,2).Value2 = "RTX3090"
$ws.Cells.Item(107,3).Value2 = 4
$ws.Cells.Item(107,4).Value2 = 128
$ws.Cells.Item(107,5).Value2 = "pannuke"
$ws.Cells.Item(107,6).Value2 = "hovernet"
$ws.Cells.Item(107,7).Value2 = "N/A"
$ws.Cells.Item(107,8).Value2 = "inference"
$ws.Cells.Item(107,9).Value2 = "TCGA-5M-AAT5-01Z-00-DX1"
$ws.Cells.Item(107,10).Value2 = "OOM"

# Row 108
$ws.Cells.Item(108,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(108,2).Value2 = "RTX3090"
$ws.Cells.Item(108,3).Value2 = 4
$ws.Cells.Item(108,4).Value2 = 128
$ws.Cells.Item(108,5).Value2 = "pannuke"
$ws.Cells.Item(108,6).Value2 = "hovernet"
$ws.Cells.Item(108,7).Value2 = "N/A"
$ws.Cells.Item(108,8).Value2 = "postproc"
$ws.Cells.Item(108,9).Value2 = "TCGA-AA-3977-01Z-00-DX1"
$ws.Cells.Item(108,10).Value2 = 1565

# Row 109
$ws.Cells.Item(109,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(109,2).Value2 = "RTX3090"
$ws.Cells.Item(109,3).Value2 = 4
$ws.Cells.Item(109,4).Value2 = 128
$ws.Cells.Item(109,5).Value2 = "pannuke"
$ws.Cells.Item(109,6).Value2 = "hovernet"
$ws.Cells.Item(109,7).Value2 = "N/A"
$ws.Cells.Item(109,8).Value2 = "postproc"
$ws.Cells.Item(109,9).Value2 = "TCGA-AA-3688-01Z-00-DX1"
$ws.Cells.Item(109,10).Value2 = "OOM"

# Row 110
$ws.Cells.Item(110,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(110,2).Value2 = "RTX3090"
$ws.Cells.Item(110,3).Value2 = 4
$ws.Cells.Item(110,4).Value2 = 128
$ws.Cells.Item(110,5).Value2 = "pannuke"
$ws.Cells.Item(110,6).Value2 = "hovernet"
$ws.Cells.Item(110,7).Value2 = "N/A"
$ws.Cells.Item(110,8).Value2 = "postproc"
$ws.Cells.Item(110,9).Value2 = "TCGA-AA-A010-01Z-00-DX1"
$ws.Cells.Item(110,10).Value2 = "OOM"

# Row 111
$ws.Cells.Item(111,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(111,2).Value2 = "RTX3090"
$ws.Cells.Item(111,3).Value2 = 4
$ws.Cells.Item(111,4).Value2 = 128
$ws.Cells.Item(111,5).Value2 = "pannuke"
$ws.Cells.Item(111,6).Value2 = "hovernet"
$ws.Cells.Item(111,7).Value2 = "N/A"
$ws.Cells.Item(111,8).Value2 = "postproc"
$ws.Cells.Item(111,9).Value2 = "TCGA-CK-4951-01Z-00-DX1"
$ws.Cells.Item(111,10).Value2 = "OOM"

# Row 112
$ws.Cells.Item(112,1).Value2 = "AMD EPYC 7302"
$ws.Cells.Item(112,2).Value2 = "RTX3090"
$ws.Cells.Item(112,3).Value2 = 4
$ws.Cells.Item(112,4).Value2 = 128
$ws.Cells.Item(112,5).Value2 = "pannuke"
$ws.Cells.Item(112,6).Value2 = "hovernet"
$ws.Cells.Item(112,7).Value2 = "N/A"
$ws.Cells.Item(112,8).Value2 = "postproc"
$ws.Cells.Item(112,9).Value2 = "TCGA-5M-AAT5-01Z-00-DX1"
$ws.Cells.Item(112,10).Value2 = "OOM"


# --- Restore view state (selection / scroll position) ---
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J9").Select()
